# "Generate Report for handoff"
#
# The 00a35824-022b-45ab-9b82-2966ccba975e.md file has been re-handed-off:
# its Status flips from "Handed back: in sync with en-US" to "Ready for
# handoff" and it gets fresh handoff timestamps on the zh-cn / de-de
# report sheets. The three tracked files are then re-sorted by file name
# on every sheet (ffffef68ae99... , ffffff799d7797... , 00a35824...),
# which is why the *other* two rows' content also shifts up by one slot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Target cell text, per worksheet / cell address.
# ---------------------------------------------------------------------
$targets = @{
    "Overview" = @{
        "A2" = "ffffef68ae99-3e08-4e05-8da5-47460aa74dd9.md"
        "B2" = "Handed back: in sync with en-US"
        "C2" = "Handed back: in sync with en-US"
        "A3" = "ffffff799d7797-f104-427c-9c24-d324d93c5151.md"
        "B3" = "Handed back: in sync with en-US"
        "C3" = "Handed back: in sync with en-US"
        "A4" = "00a35824-022b-45ab-9b82-2966ccba975e.md"
        "B4" = "Ready for handoff"
        "C4" = "Ready for handoff"
    }
    "zh-cn" = @{
        "A2" = "ffffef68ae99-3e08-4e05-8da5-47460aa74dd9.md"
        "C2" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.306683e0ae1bc6e6d5e469209d2ad1cee402f44d.zh-cn.xlf"
        "D2" = "2016-01-18 11:26:52"
        "E2" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.md"
        "F2" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.306683e0ae1bc6e6d5e469209d2ad1cee402f44d.zh-cn.xlf"
        "G2" = "2016-01-18 11:27:33"
        "A3" = "ffffff799d7797-f104-427c-9c24-d324d93c5151.md"
        "C3" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.306683e0ae1bc6e6d5e469209d2ad1cee402f44d.zh-cn.xlf"
        "D3" = "2016-01-18 11:26:52"
        "E3" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.md"
        "F3" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.306683e0ae1bc6e6d5e469209d2ad1cee402f44d.zh-cn.xlf"
        "G3" = "2016-01-18 11:27:33"
        "A4" = "00a35824-022b-45ab-9b82-2966ccba975e.md"
        "B4" = "Ready for handoff"
        "C4" = "00a35824-022b-45ab-9b82-2966ccba975e.0b4fde6d1f068c6dd722e1b2eedbf47a3c2b7b59.zh-cn.xlf"
        "D4" = "2016-01-18 11:32:01"
        "E4" = "00a35824-022b-45ab-9b82-2966ccba975e.md"
        "F4" = "00a35824-022b-45ab-9b82-2966ccba975e.0b4fde6d1f068c6dd722e1b2eedbf47a3c2b7b59.zh-cn.xlf"
        "G4" = "2016-01-18 11:31:03"
    }
    "de-de" = @{
        "A2" = "ffffef68ae99-3e08-4e05-8da5-47460aa74dd9.md"
        "C2" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.306683e0ae1bc6e6d5e469209d2ad1cee402f44d.de-de.xlf"
        "D2" = "2016-01-18 11:27:02"
        "E2" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.md"
        "F2" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.306683e0ae1bc6e6d5e469209d2ad1cee402f44d.de-de.xlf"
        "G2" = "2016-01-18 11:27:49"
        "A3" = "ffffff799d7797-f104-427c-9c24-d324d93c5151.md"
        "C3" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.306683e0ae1bc6e6d5e469209d2ad1cee402f44d.de-de.xlf"
        "D3" = "2016-01-18 11:27:02"
        "E3" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.md"
        "F3" = "ee3d5f9e-b5ed-409a-87a4-7e48ee1f7f72.306683e0ae1bc6e6d5e469209d2ad1cee402f44d.de-de.xlf"
        "G3" = "2016-01-18 11:27:49"
        "A4" = "00a35824-022b-45ab-9b82-2966ccba975e.md"
        "B4" = "Ready for handoff"
        "C4" = "00a35824-022b-45ab-9b82-2966ccba975e.0b4fde6d1f068c6dd722e1b2eedbf47a3c2b7b59.de-de.xlf"
        "D4" = "2016-01-18 11:32:12"
        "E4" = "00a35824-022b-45ab-9b82-2966ccba975e.md"
        "F4" = "00a35824-022b-45ab-9b82-2966ccba975e.0b4fde6d1f068c6dd722e1b2eedbf47a3c2b7b59.de-de.xlf"
        "G4" = "2016-01-18 11:31:19"
    }
}

foreach ($sheetName in $targets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellTargets = $targets[$sheetName]

    # Write the new literal cell text for every changed cell.
    foreach ($addr in $cellTargets.Keys) {
        $ws.Range($addr).Value2 = $cellTargets[$addr]
    }

    # Hyperlinks keep their original target URL (the relationship file
    # didn't change) but the visible display text has to track whatever
    # we just wrote into column A/C/E/F. Update in place, live off the
    # enumerator -- indexed access re-creates a disconnected wrapper that
    # can't be written back.
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address($false, $false)
        if ($cellTargets.ContainsKey($addr)) {
            $h.TextToDisplay = $cellTargets[$addr]
        }
    }
}
